$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.890.36'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '2.733.05'
$ws.Range('E3').Value = '  +3.48%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.06'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.22'
$ws.Range('E6').Value = '  +6.45%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('D9').Value = '2.731.45'
$ws.Range('E9').Value = '  +3.46%  '
$ws.Range('E10').Value = '  +2.93%  '
$ws.Range('E11').Value = '  +4.77%  '
$ws.Range('E12').Value = '  +1.67%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.68'
$ws.Range('E14').Value = '  +3.10%  '
$ws.Range('D15').Value = '3.232.26'
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('D17').Value = '68.812.70'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').Value = '2.777.00'
$ws.Range('E18').Value = '  +6.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.85'
$ws.Range('E19').Value = '  +4.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '373.49'
$ws.Range('E20').Value = '  +3.83%  '
$ws.Range('E21').Value = '  +5.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.54'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.96'
$ws.Range('E23').Value = '  +4.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  +3.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.73'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.98'
$ws.Range('E27').Value = '  +2.33%  '
$ws.Range('D28').Value = '2.877.22'
$ws.Range('E28').Value = '  +3.76%  '
$ws.Range('E29').Value = '  +3.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '592.37'
$ws.Range('E30').Value = '  +5.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.28'
$ws.Range('E32').Value = '  +3.87%  '
$ws.Range('E33').Value = '  +5.19%  '
$ws.Range('E34').Value = '  +6.26%  '
$ws.Range('E35').Value = '  +4.67%  '
$ws.Range('E36').Value = '  +5.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '161.13'
$ws.Range('E38').Value = '  +1.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.92'
$ws.Range('E39').Value = '  +1.26%  '
$ws.Range('E40').Value = '  +3.34%  '
$ws.Range('E41').Value = '  +3.47%  '
$ws.Range('E42').Value = '  +3.94%  '
$ws.Range('E43').Value = '  +3.29%  '
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').Value = '0.0₆0314'
$ws.Range('E46').Value = '  -2.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '41.01'
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '157.07'
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('E49').Value = '  +6.63%  '
$ws.Range('E50').Value = '  +7.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.606'
$ws.Range('E51').Value = '  +7.49%  '
